$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Kit")

# B25: R39,R40,R54 -> R54
$ws.Range("B25").Value = "R54"

# B26: ...R50,R51,R57,R58 -> add R39,R40 before R50
$ws.Range("B26").Value = "R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R40,R50,R51,R57,R58"

# B27: R9,R12,R15,R18,R26,R28,R33,R34 -> R9,R12,R15,R18
$ws.Range("B27").Value = "R9,R12,R15,R18"

# B29: R1,R3,R59 -> R1,R3,R26,R28,R33,R34,R59 ; also turn A29 into a formula like the others
$ws.Range("B29").Value = "R1,R3,R26,R28,R33,R34,R59"
$ws.Range("A29").Formula = '=LEN(B29)-LEN(SUBSTITUTE(B29,",",""))+1'

$ws.Rows.Item(26).RowHeight = 27
